$wb = $excel.ActiveWorkbook
$wsMain  = $wb.Worksheets.Item("Main")
$wsFLRbI = $wb.Worksheets.Item("FLRbI")

# --- Text / content updates -------------------------------------------------
# Replace the old "not attempted ... California" note with the new explanatory
# text, and add a second note about the California model inheriting values.
$wsMain.Range("A3").Value = "The model structure calculates leakage but does not integrated in policy impcts."
$wsMain.Range("A4").Value = "California model retains values for this variable inherited from the US EPS."

# --- Strip the old placeholder/formatting rows below row 4 -----------------
# The sheet used to carry a long tail of empty, pre-formatted rows (yellow
# highlight blocks, a hyperlink-styled cell, bold placeholder cells, etc.)
# down to row 45. Remove all of that.
$wsMain.Rows("5:45").Delete()

# Recreate the three remaining bold placeholder cells at their new positions.
$wsMain.Range("A12").Font.Bold = $true
$wsMain.Range("A22").Font.Bold = $true
$wsMain.Range("A29").Font.Bold = $true

# Restore row 3's slightly shorter custom height.
$wsMain.Rows(3).RowHeight = 14.45

# --- Calculation settings ---------------------------------------------------
# Turn off iterative calculation (was on for the old circular-reference
# workaround).
$excel.Iteration = $false

# --- Active sheet / selection ----------------------------------------------
# Main becomes the active/visible sheet (was FLRbI), with A5 selected.
$wsMain.Activate()
$wsMain.Range("A5").Select()

# FLRbI keeps its own selection on B5 (unchanged) but is no longer the active tab.
$wsFLRbI.Range("B5").Select()
$wsMain.Activate()
